# Weekly crime-stat update: new week of data rolled in (Volume 32, Number 20;
# report window 5/12/2025 - 5/18/2025), refreshing the Week-to-Date / 28-Day /
# Year-to-Date / 2-Year figures and their %-change columns on the CompStat sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text (shared strings made of multiple same-format runs; a
#     plain concatenated string reproduces the same rendered text) ---
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

$cellData = @(
    @{cell="C14"; value=1; fmt="int"},
    @{cell="F14"; value=1; fmt="int"},
    @{cell="I14"; value=1; fmt="int"},
    @{cell="K14"; value=0; fmt="pct"},
    @{cell="L14"; value=-75; fmt="pct"},
    @{cell="M14"; value=-80; fmt="pct"},
    @{cell="N14"; value=-88.888888888888; fmt="pct"},
    @{cell="C15"; value=2; fmt="int"},
    @{cell="D15"; value=1; fmt="int"},
    @{cell="E15"; value=100; fmt="pct"},
    @{cell="F15"; value=5; fmt="int"},
    @{cell="G15"; value=3; fmt="int"},
    @{cell="H15"; value=66.666666666666; fmt="pct"},
    @{cell="I15"; value=15; fmt="int"},
    @{cell="J15"; value=17; fmt="int"},
    @{cell="K15"; value=-11.764705882352; fmt="pct"},
    @{cell="L15"; value=-11.764705882352; fmt="pct"},
    @{cell="M15"; value=200; fmt="pct"},
    @{cell="N15"; value=-37.5; fmt="pct"},
    @{cell="C16"; value=6; fmt="int"},
    @{cell="E16"; value=-45.454545454545; fmt="pct"},
    @{cell="F16"; value=43; fmt="int"},
    @{cell="G16"; value=45; fmt="int"},
    @{cell="H16"; value=-4.444444444444; fmt="pct"},
    @{cell="I16"; value=169; fmt="int"},
    @{cell="J16"; value=161; fmt="int"},
    @{cell="K16"; value=4.968944099378; fmt="pct"},
    @{cell="L16"; value=19.858156028368; fmt="pct"},
    @{cell="M16"; value=79.787234042553; fmt="pct"},
    @{cell="N16"; value=-61.415525114155; fmt="pct"},
    @{cell="C17"; value=17; fmt="int"},
    @{cell="D17"; value=13; fmt="int"},
    @{cell="E17"; value=30.76923076923; fmt="pct"},
    @{cell="F17"; value=73; fmt="int"},
    @{cell="G17"; value=65; fmt="int"},
    @{cell="H17"; value=12.307692307692; fmt="pct"},
    @{cell="I17"; value=264; fmt="int"},
    @{cell="J17"; value=274; fmt="int"},
    @{cell="K17"; value=-3.649635036496; fmt="pct"},
    @{cell="L17"; value=14.285714285714; fmt="pct"},
    @{cell="M17"; value=169.387755102041; fmt="pct"},
    @{cell="N17"; value=-23.478260869565; fmt="pct"},
    @{cell="C18"; value=3; fmt="int"},
    @{cell="D18"; value=4; fmt="int"},
    @{cell="E18"; value=-25; fmt="pct"},
    @{cell="F18"; value=16; fmt="int"},
    @{cell="H18"; value=-30.434782608695; fmt="pct"},
    @{cell="I18"; value=95; fmt="int"},
    @{cell="J18"; value=107; fmt="int"},
    @{cell="K18"; value=-11.214953271028; fmt="pct"},
    @{cell="L18"; value=26.666666666666; fmt="pct"},
    @{cell="M18"; value=97.916666666666; fmt="pct"},
    @{cell="N18"; value=-75.949367088607; fmt="pct"},
    @{cell="C19"; value=13; fmt="int"},
    @{cell="D19"; value=20; fmt="int"},
    @{cell="E19"; value=-35; fmt="pct"},
    @{cell="F19"; value=48; fmt="int"},
    @{cell="G19"; value=51; fmt="int"},
    @{cell="H19"; value=-5.882352941176; fmt="pct"},
    @{cell="I19"; value=200; fmt="int"},
    @{cell="J19"; value=230; fmt="int"},
    @{cell="K19"; value=-13.043478260869; fmt="pct"},
    @{cell="L19"; value=21.951219512195; fmt="pct"},
    @{cell="M19"; value=138.095238095238; fmt="pct"},
    @{cell="N19"; value=65.289256198347; fmt="pct"},
    @{cell="C20"; value=2; fmt="int"},
    @{cell="D20"; value=2; fmt="int"},
    @{cell="E20"; value=0; fmt="pct"},
    @{cell="F20"; value=17; fmt="int"},
    @{cell="G20"; value=21; fmt="int"},
    @{cell="H20"; value=-19.047619047619; fmt="pct"},
    @{cell="I20"; value=96; fmt="int"},
    @{cell="J20"; value=91; fmt="int"},
    @{cell="K20"; value=5.494505494505; fmt="pct"},
    @{cell="L20"; value=-45.762711864406; fmt="pct"},
    @{cell="M20"; value=113.333333333333; fmt="pct"},
    @{cell="N20"; value=-54.502369668246; fmt="pct"},
    @{cell="C21"; value=44; fmt="none"},
    @{cell="D21"; value=51; fmt="none"},
    @{cell="E21"; value=-13.725490196078; fmt="none"},
    @{cell="F21"; value=203; fmt="none"},
    @{cell="G21"; value=208; fmt="none"},
    @{cell="H21"; value=-2.403846153846; fmt="none"},
    @{cell="I21"; value=840; fmt="none"},
    @{cell="J21"; value=881; fmt="none"},
    @{cell="K21"; value=-4.653802497162; fmt="none"},
    @{cell="L21"; value=3.831891223733; fmt="none"},
    @{cell="M21"; value=121.635883905013; fmt="none"},
    @{cell="N21"; value=-45.560596241088; fmt="none"},
    @{cell="M22"; value=-77.777777777777; fmt="pct"},
    @{cell="C23"; value=4; fmt="int"},
    @{cell="D23"; value=7; fmt="int"},
    @{cell="E23"; value=-42.857142857142; fmt="pct"},
    @{cell="F23"; value=26; fmt="int"},
    @{cell="G23"; value=28; fmt="int"},
    @{cell="H23"; value=-7.142857142857; fmt="pct"},
    @{cell="I23"; value=110; fmt="int"},
    @{cell="J23"; value=147; fmt="int"},
    @{cell="K23"; value=-25.17006802721; fmt="pct"},
    @{cell="L23"; value=-28.571428571428; fmt="pct"},
    @{cell="M23"; value=59.420289855072; fmt="pct"},
    @{cell="C24"; value=15; fmt="int"},
    @{cell="D24"; value=24; fmt="int"},
    @{cell="E24"; value=-37.5; fmt="pct"},
    @{cell="F24"; value=98; fmt="int"},
    @{cell="G24"; value=82; fmt="int"},
    @{cell="H24"; value=19.512195121951; fmt="pct"},
    @{cell="I24"; value=451; fmt="int"},
    @{cell="J24"; value=419; fmt="int"},
    @{cell="K24"; value=7.637231503579; fmt="pct"},
    @{cell="L24"; value=3.67816091954; fmt="pct"},
    @{cell="M24"; value=61.648745519713; fmt="pct"},
    @{cell="C25"; value=8; fmt="int"},
    @{cell="E25"; value=700; fmt="pct"},
    @{cell="G25"; value=8; fmt="int"},
    @{cell="H25"; value=262.5; fmt="pct"},
    @{cell="I25"; value=94; fmt="int"},
    @{cell="J25"; value=60; fmt="int"},
    @{cell="K25"; value=56.666666666666; fmt="pct"},
    @{cell="L25"; value=1.075268817204; fmt="pct"},
    @{cell="C26"; value=16; fmt="int"},
    @{cell="D26"; value=13; fmt="int"},
    @{cell="E26"; value=23.076923076923; fmt="pct"},
    @{cell="F26"; value=66; fmt="int"},
    @{cell="G26"; value=66; fmt="int"},
    @{cell="I26"; value=326; fmt="int"},
    @{cell="J26"; value=320; fmt="int"},
    @{cell="K26"; value=1.875; fmt="pct"},
    @{cell="L26"; value=-24.884792626728; fmt="pct"},
    @{cell="M26"; value=4.823151125401; fmt="pct"},
    @{cell="C27"; value=2; fmt="int"},
    @{cell="E27"; value=100; fmt="pct"},
    @{cell="F27"; value=5; fmt="int"},
    @{cell="G27"; value=4; fmt="int"},
    @{cell="H27"; value=25; fmt="pct"},
    @{cell="I27"; value=16; fmt="int"},
    @{cell="J27"; value=24; fmt="int"},
    @{cell="K27"; value=-33.333333333333; fmt="pct"},
    @{cell="L27"; value=-36; fmt="pct"},
    @{cell="C28"; value=1; fmt="int"},
    @{cell="D28"; value=2; fmt="int"},
    @{cell="E28"; value=-50; fmt="pct"},
    @{cell="G28"; value=14; fmt="int"},
    @{cell="H28"; value=-71.428571428571; fmt="pct"},
    @{cell="I28"; value=20; fmt="int"},
    @{cell="J28"; value=35; fmt="int"},
    @{cell="K28"; value=-42.857142857142; fmt="pct"},
    @{cell="L28"; value=-51.219512195122; fmt="pct"},
    @{cell="C29"; value=3; fmt="int"},
    @{cell="D29"; value=2; fmt="int"},
    @{cell="E29"; value=50; fmt="pct"},
    @{cell="F29"; value=3; fmt="int"},
    @{cell="G29"; value=5; fmt="int"},
    @{cell="H29"; value=-40; fmt="pct"},
    @{cell="I29"; value=10; fmt="int"},
    @{cell="J29"; value=9; fmt="int"},
    @{cell="K29"; value=11.111111111111; fmt="pct"},
    @{cell="L29"; value=0; fmt="pct"},
    @{cell="M29"; value=-47.368421052631; fmt="pct"},
    @{cell="N29"; value=-66.666666666666; fmt="pct"},
    @{cell="C30"; value=2; fmt="int"},
    @{cell="D30"; value=2; fmt="int"},
    @{cell="E30"; value=0; fmt="pct"},
    @{cell="F30"; value=2; fmt="int"},
    @{cell="G30"; value=5; fmt="int"},
    @{cell="H30"; value=-60; fmt="pct"},
    @{cell="I30"; value=8; fmt="int"},
    @{cell="J30"; value=9; fmt="int"},
    @{cell="K30"; value=-11.111111111111; fmt="pct"},
    @{cell="L30"; value=-20; fmt="pct"},
    @{cell="M30"; value=-46.666666666666; fmt="pct"},
    @{cell="N30"; value=-73.333333333333; fmt="pct"},
)

$intFormat = "#,##0"
$pctFormat = '#,##0.0;"-"#,##0.0'

foreach ($item in $cellData) {
    $rng = $ws.Range($item.cell)
    $rng.Value = $item.value
    if ($item.fmt -eq "int") {
        $rng.NumberFormat = $intFormat
    } elseif ($item.fmt -eq "pct") {
        $rng.NumberFormat = $pctFormat
    }
}

# Column H narrows back down now that this week's %-change figures are no
# longer the widest entries in the column (matches columns F/G/I/J).
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
